$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1950
$ws.Range("I17").Value = 1950
$ws.Range("K17").Value = 5850
$ws.Range("M17").Value = -5682
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766
$ws.Range("H38").Value = 1072.5
$ws.Range("J38").Value = 1999.6666
$ws.Range("L38").Value = 5998.9998
$ws.Range("N38").Value = -6742.9998
$ws.Range("H40").Value = 6346.5386
$ws.Range("I40").Value = 4000.2
$ws.Range("J40").Value = 7813
$ws.Range("K40").Value = 4000.2
$ws.Range("L40").Value = 7813
$ws.Range("M40").Value = -3825.2
$ws.Range("N40").Value = -8163
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = $null
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = $null
$ws.Range("N71").Value = 0
$ws.Range("H98").Value = 942.4286
$ws.Range("I98").Value = 942.4286
$ws.Range("K98").Value = 942.4286
$ws.Range("M98").Value = 555.5714
$ws.Range("H107").Value = 1128.8125
$ws.Range("I107").Value = 747.3
$ws.Range("J107").Value = 1764.6666
$ws.Range("K107").Value = 747.3
$ws.Range("L107").Value = 1764.6666
$ws.Range("M107").Value = 1172.7
$ws.Range("N107").Value = -5604.6666
$ws.Range("H122").Value = 942.4286
$ws.Range("I122").Value = 942.4286
$ws.Range("K122").Value = 2827.2858
$ws.Range("M122").Value = -377.2857999999997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4216.2173
$ws.Range("I32").Value = 2662.7896
$ws.Range("K32").Value = 2662.7896
$ws.Range("M32").Value = -2375.7896
$ws.Range("H61").Value = 2913.6
$ws.Range("I61").Value = 2715.111
$ws.Range("K61").Value = 2715.111
$ws.Range("M61").Value = -2503.111
$ws.Range("H136").Value = 2913.6
$ws.Range("I136").Value = 2715.111
$ws.Range("K136").Value = 8145.333
$ws.Range("M136").Value = -5595.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4799.7144
$ws.Range("I86").Value = 1899.5
$ws.Range("J86").Value = 5959.8
$ws.Range("K86").Value = 1899.5
$ws.Range("L86").Value = 5959.8
$ws.Range("M86").Value = -776.5
$ws.Range("N86").Value = -8205.799999999999
$ws.Range("H89").Value = 4799.7144
$ws.Range("I89").Value = 1899.5
$ws.Range("J89").Value = 5959.8
$ws.Range("K89").Value = 9497.5
$ws.Range("L89").Value = 29799
$ws.Range("M89").Value = -3881.5
$ws.Range("N89").Value = -41031

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 229.5
$ws.Range("I7").Value = 229.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 229.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -116.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 871.75
$ws.Range("I8").Value = 871.75
$ws.Range("K8").Value = 2615.25
$ws.Range("M8").Value = -2476.25
$ws.Range("H12").Value = 43.833332
$ws.Range("I12").Value = 53.5
$ws.Range("K12").Value = 160.5
$ws.Range("M12").Value = 12.5
$ws.Range("H80").Value = 6124.875
$ws.Range("J80").Value = 6124.875
$ws.Range("L80").Value = 18374.625
$ws.Range("N80").Value = -20246.625
$ws.Range("H83").Value = 6124.875
$ws.Range("J83").Value = 6124.875
$ws.Range("L83").Value = 55123.875
$ws.Range("N83").Value = -64483.875
$ws.Range("H98").Value = 664.6667
$ws.Range("I98").Value = 664.6667
$ws.Range("K98").Value = 1994.0001
$ws.Range("M98").Value = -496.0001
$ws.Range("H113").Value = 1196
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = $null
$ws.Range("H122").Value = 4405
$ws.Range("J122").Value = 4405
$ws.Range("L122").Value = 39645
$ws.Range("N122").Value = -44545
$ws.Range("H132").Value = 2400.2856
$ws.Range("I132").Value = 2360.4
$ws.Range("K132").Value = 21243.6
$ws.Range("M132").Value = -18713.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3245.6155
$ws.Range("I102").Value = 2949.5
$ws.Range("J102").Value = 4232.6665
$ws.Range("K102").Value = 2949.5
$ws.Range("L102").Value = 4232.6665
$ws.Range("M102").Value = -1327.5
$ws.Range("N102").Value = -7476.6665
$ws.Range("H122").Value = 8828.143
$ws.Range("I122").Value = 8828.143
$ws.Range("K122").Value = 26484.429
$ws.Range("M122").Value = -24034.429
$ws.Range("H132").Value = 3689.5334
$ws.Range("I132").Value = 3746.5386
$ws.Range("K132").Value = 11239.6158
$ws.Range("M132").Value = -8709.6158

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 30004
$ws.Range("I3").Value = 30004
$ws.Range("K3").Value = 30004
$ws.Range("M3").Value = -29892
$ws.Range("H15").Value = 30004
$ws.Range("I15").Value = 30004
$ws.Range("K15").Value = 30004
$ws.Range("M15").Value = -29834
$ws.Range("H98").Value = 40000
$ws.Range("J98").Value = 40000
$ws.Range("L98").Value = 40000
$ws.Range("N98").Value = -45990
$ws.Range("H132").Value = 19035
$ws.Range("I132").Value = 22221.889
$ws.Range("J132").Value = 4694
$ws.Range("K132").Value = 66665.667
$ws.Range("L132").Value = 14082
$ws.Range("M132").Value = -64135.667
$ws.Range("N132").Value = -19142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 30006
$ws.Range("I18").Value = 30006
$ws.Range("K18").Value = 30006
$ws.Range("M18").Value = -29833
$ws.Range("H104").Value = 10123.333
$ws.Range("J104").Value = 10123.333
$ws.Range("L104").Value = 10123.333
$ws.Range("N104").Value = -17111.333
$ws.Range("H132").Value = 3352.7
$ws.Range("I132").Value = 3597.4443
$ws.Range("J132").Value = 1150
$ws.Range("K132").Value = 10792.3329
$ws.Range("L132").Value = 3450
$ws.Range("M132").Value = -8262.332900000001
$ws.Range("N132").Value = -8510
